$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the new metadata sheet ---
$datos = $wb.Worksheets.Item(1)
$datos.Name = "datos"

$metadatos = $wb.Worksheets.Add($null, $datos)
$metadatos.Name = "metadatos"

# --- "datos" sheet: drop the bold/bordered header style, keep values ---
$datos.Range("A1:C1").ClearFormats()
$datos.Range("G7").Select()

# --- "metadatos" sheet content ---
# Row 1 - column headers
$metadatos.Range("A1").Value = "Variables"
$metadatos.Range("B1").Value = "Descripción"
$metadatos.Range("C1").Value = "Fuente"
$metadatos.Range("D1").Value = "Fecha_de_extracción"

# Row 2 - anno
$metadatos.Range("A2").Value = "anno"
$metadatos.Range("B2").Value = "Año"
$metadatos.Range("C2").Value = "…"
$metadatos.Range("D2").Value = 45715

# Row 3 - codmpio
$metadatos.Range("A3").Value = "codmpio"
$metadatos.Range("B3").Value = "Código del municipio"
$metadatos.Range("C3").Value = "…"
$metadatos.Range("D3").Value = 45715

# Row 4 - pobreza_monetaria
$metadatos.Range("A4").Value = "pobreza_monetaria"
$metadatos.Range("B4").Value = "Porcentaje"
$metadatos.Range("C4").Value = "Departamento Nacional de Planeación (DNP)-TERRIDATA"
$metadatos.Range("D4").Value = 45715

# Formatting: black Calibri 11 text for labels, black date format for dates
$metadatos.Range("D2:D4").Font.Color = 0
$metadatos.Range("D2:D4").NumberFormat = "mm-dd-yy"

$metadatos.Range("A1:D1").Font.Color = 0
$metadatos.Range("A2:C3").Font.Color = 0
$metadatos.Range("B4:C4").Font.Color = 0

$metadatos.Range("C4").Select()
$metadatos.Activate()

$wb.Save()
